$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the 95b95732-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 22:16:35"

# --- zh-cn sheet: row 3 is the 95b95732-... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 22:16:29"

# --- de-de sheet: row 3 is the 95b95732-... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 22:16:35"
